# "Doing Updates for Financials" — refresh the PRTHU yearly financials
# figures in column D (and, for a handful of rows that were previously
# all-zero placeholder rows, recast the now-inapplicable E:J figures as
# "NA" to match the existing NA convention used elsewhere in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRTHU")

# --- Income Statement -------------------------------------------------
# Income After Tax
$ws.Range("D21").Value = 14700
$ws.Range("E21:J21").Value = "NA"

# --- Balance Sheet ------------------------------------------------------
# Cash And Cash Equivalents
$ws.Range("D41").Value = 28100
# Net Receivables
$ws.Range("D43").Value = 51100
$ws.Range("E43:J43").Value = "NA"
# Other Current Assets
$ws.Range("D45").Value = 26900
# Total Current Assets
$ws.Range("D46").Value = 106000
# Long Term Investments
$ws.Range("D47").Value = 3800
$ws.Range("E47:J47").Value = "NA"
# Property Plant and Equipment
$ws.Range("D48").Value = 11900
$ws.Range("E48:J48").Value = "NA"
# Goodwill
$ws.Range("D49").Value = 143600
$ws.Range("E49:J49").Value = "NA"
# Other Assets
$ws.Range("D52").Value = 56500
# Total Assets
$ws.Range("D54").Value = 266700
# Accounts Payable
$ws.Range("D57").Value = 18600
# Short/Current Long Term Debt
$ws.Range("D58").Value = 9100
# Other Current Liabilities
$ws.Range("D59").Value = 38800
# Total Current Liabilities
$ws.Range("D60").Value = 66500
# Long Term Debt
$ws.Range("D61").Value = 275600
# Other Liabilities
$ws.Range("D62").Value = 15800
# Total Liabilities
$ws.Range("D66").Value = 356900
# Common Stock
$ws.Range("D72").Value = -90200
# Total Stockholder Equity
$ws.Range("D76").Value = -90200

# --- Cash Flow Statement -------------------------------------------------
# Depreciation
$ws.Range("D83").Value = 14700
$ws.Range("E83:J83").Value = "NA"
# Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 36900
# Capital Expenditures
$ws.Range("D91").Value = -6600
$ws.Range("E91:J91").Value = "NA"
# Other Cashflows from Investing Activities
$ws.Range("D94").Value = -9000
# Dividends Paid
$ws.Range("D96").Value = -3400
# Net Borrowings
$ws.Range("D100").Value = -25400
# Total Cash Flows From Financing Activities
$ws.Range("J101").Value = "NA"
# Change In Cash and Cash Equivalents
$ws.Range("D102").Value = 2500
